$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B2 username
$ws.Range("B2").Value = "saalim192"

# Update row 3 with new data
$ws.Range("A3").Value = "'957228261610496001"
$ws.Range("B3").Value = "faizanr7"
$ws.Range("C3").Value = "Twitter"
$ws.Range("D3").Value = 1

# Delete rows 4 through 7 which are no longer needed
$ws.Range("A4:D7").EntireRow.Delete()
